$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.954.32"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "1.632.96"
$ws.Range("E3").Value = "  -2.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.38"
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5100"
$ws.Range("E6").Value = "  -1.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.006"
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2541"
$ws.Range("E8").Value = "  -3.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06139"
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.32"
$ws.Range("E10").Value = "  -3.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07558"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "1.646.29"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.332"
$ws.Range("E13").Value = "  -1.97%  "
$ws.Range("D14").Value = "1.861.96"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5318"
$ws.Range("E15").Value = "  -4.91%  "
$ws.Range("D16").Value = "0.0₅7937"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.03"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").Value = "25.961.97"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.604"
$ws.Range("E20").Value = "  -3.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "184.05"
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.914"
$ws.Range("E22").Value = "  -4.22%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.006"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.046"
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.64"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1186"
$ws.Range("E26").Value = "  -4.56%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.270"
$ws.Range("E27").Value = "  -3.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.38"
$ws.Range("E28").Value = "  -3.70%  "
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05999"
$ws.Range("E30").Value = "  -3.92%  "
$ws.Range("E31").Value = "  -2.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.375"
$ws.Range("E32").Value = "  -2.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.338"
$ws.Range("E33").Value = "  -2.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.606"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9603"
$ws.Range("E35").Value = "  -3.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.381"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.714"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5783"
$ws.Range("E38").Value = "  -4.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01573"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("D40").Value = "1.073.06"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.774"
$ws.Range("E41").Value = "  -5.76%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8401"
$ws.Range("E43").Value = "  -2.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.61"
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("D46").Value = "0.0₈108"
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9973"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.94"
$ws.Range("E48").Value = "  -3.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.923"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05202"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4232"
$ws.Range("E51").Value = "  -0.38%  "
